# "added new version of completed powerapp gym"
#
# 1. Products sheet: the Revenue row's Photo hyperlink (G4) pointed at a
#    thumbnail image (03_Revenue_tn.jpg); it now points at the full-size,
#    final image (03_Revenue.jpg).
# 2. Subscriptions sheet/table: drop the leftover Power Apps plumbing
#    column (__PowerAppsId__ / "7GvcffR1DnU") that shouldn't ship with the
#    finished app.
# 3. Refresh the view state: Subscriptions is now the active/selected tab,
#    with I10 selected there, while Products remembers G4 selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Products")
$ws2 = $wb.Worksheets.Item("Subscriptions")

# --- 1. Fix the Revenue product's photo URL -----------------------------
$ws1.Range("G4").Value = "https://raw.githubusercontent.com/xpandit/powerplatform-powergym/master/Assets/03_Revenue.jpg"

# --- 2. Remove the __PowerAppsId__ column from the Subscription table ---
$lo2 = $ws2.ListObjects.Item("Subscription")
$lo2.ListColumns.Item("__PowerAppsId__").Delete()
# Clear out the now-empty column so the sheet's used range shrinks too.
$ws2.Columns.Item(4).Delete()

# --- 3. Update selection / active-sheet state ----------------------------
$ws1.Activate()
$ws1.Range("G4").Select()

$ws2.Activate()
$ws2.Range("I10").Select()
